$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.859.21'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '2.300.60'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.22%  '
$ws.Range('E7').Value = '  +1.17%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +4.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.81'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.31%  '
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +12.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.34%  '
$ws.Range('D15').Value = '2.661.39'
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').Value = '2.289.82'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.801'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.97%  '
$ws.Range('D18').Value = '42.757.99'
$ws.Range('E18').Value = '  +1.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.21%  '
$ws.Range('D21').Value = '0.0₃0898'
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E24').Value = '  +12.88%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.52'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.30'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.12'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.61'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.29'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.56%  '
$ws.Range('E36').Value = '  +3.33%  '
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('E38').Value = '  +3.40%  '
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('E40').Value = '  +3.03%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').Value = '1.977.39'
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('E43').Value = '  +3.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.14'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.60'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.78%  '
$ws.Range('D49').Value = '2.528.12'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('E50').Value = '  +3.43%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.51'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.84%  '
